$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 10 ("ENTAILMENT") - fill in the 10 query run values
$ws.Range("B10").Value = 185620
$ws.Range("C10").Value = 184710
$ws.Range("D10").Value = 185580
$ws.Range("E10").Value = 182315
$ws.Range("F10").Value = 183207
$ws.Range("G10").Value = 183886
$ws.Range("H10").Value = 181302
$ws.Range("I10").Value = 181158
$ws.Range("J10").Value = 182356
$ws.Range("K10").Value = 183611

# Apply wrap-text styling to the newly-populated run cells (B10:K10)
$ws.Range("B10:K10").WrapText = $true

# Average formula in M10, matching the pattern used by the other rows
$ws.Range("M10").Formula = "=AVERAGE(B10:K10)"

# Update the active selection to N12
$ws.Range("N12").Select() | Out-Null
